# Swap the presentation's theme from "Integral" to the default "Office Theme".
#
# The authored edit replaces the colour scheme that theme1.xml (the theme
# bound to the slide master / the whole deck) carries: it goes from the
# green/olive "Integral" palette to the stock blue "Office" palette. (The
# OOXML diff also shows the old "Integral" palette being preserved as
# theme2.xml - the notes master's theme - but that part of the package
# isn't reachable for editing through the PowerPoint object model, so we
# apply the part of the change that *is* exposed: the design/theme colours
# of the presentation itself.)

$p = $ppt.ActivePresentation

# Helper: build the packed BGR integer PowerPoint's RGBColor.RGB expects
# from a "RRGGBB" hex string.
function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Standard Office theme palette, in ThemeColorScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgb $officeTheme[$i - 1]
}
